$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: 번호 -> 랭킹, 성적 -> 별점, clear 상위권 (D1)
$ws.Range("A1").Value = "랭킹"
$ws.Range("C1").Value = "별점"
$ws.Range("D1").ClearContents()
